# dsa two pointers practice
# Add a new row (125. Valid Palindrome) to the LeetCode tracking table on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Expand the existing table (ListObject) by one row. This automatically
# grows the table ref/autoFilter and the worksheet dimension.
$tbl = $ws.ListObjects.Item(1)
$newListRow = $tbl.ListRows.Add()

$rowNum = $newListRow.Range.Row

$question = "125. Valid Palindrome"
$difficulty = "Easy"
$pattern = "Two Pointers"
$notes = "Remember s = s.toLowerCase(). Alphanumeric is a-z and 0-9 and can be computed with (c >= 'a' && c <= 'z') || (c >= '0' && c <= '9'). Make a helper function. In the while loop to skip over blank characters, use continue; There is also the Character.isLetterOrDigit() if you can remember it. Watch the while loop conditions."
$link = "https://leetcode.com/problems/valid-palindrome/solutions/3165353/beats-96-9-well-explained-code-in-java/ "

$ws.Cells.Item($rowNum, 1).Value = $question

# Add the hyperlink for the Link column (column E) first, matching the
# style of the rest of the "Link" column.
$linkCell = $ws.Cells.Item($rowNum, 5)
$ws.Hyperlinks.Add($linkCell, $link)

$ws.Cells.Item($rowNum, 2).Value = $difficulty
$ws.Cells.Item($rowNum, 3).Value = $pattern
$ws.Cells.Item($rowNum, 4).Value = $notes

# Match formatting conventions used by the other rows of the table:
#  - Difficulty column uses a colored fill depending on difficulty (Easy = green, style from row above)
#  - Link column uses the Hyperlink cell style (style from row above)
$ws.Range("B" + ($rowNum - 1)).Copy()
$ws.Cells.Item($rowNum, 2).PasteSpecial(-4122)

$ws.Range("E" + ($rowNum - 1)).Copy()
$ws.Cells.Item($rowNum, 5).PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Move the active selection the way it ended up after the edit.
$ws.Range("D57").Select()
